# Auto-optimize exam scheduling: dynamically adjusts exams per slot (1-4)
# to guarantee all courses are scheduled within date range.
#
# Updates the timetable grids on Section_A and Section_B with the
# re-balanced exam slot assignments.

$wb = $excel.ActiveWorkbook

# --- Section_A (sheet1) ---------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "CS304"
$wsA.Range("C2").Value = "CS304"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "CS303"
$wsA.Range("F2").Value = "CS309"

$wsA.Range("B3").Value = "Free"
$wsA.Range("C3").Value = "CS309"
$wsA.Range("D3").Value = "Free"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "Free"

$wsA.Range("B5").Value = "CS303"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "Free"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "Free"

$wsA.Range("B6").Value = "Free"
$wsA.Range("C6").Value = "CS309 (Tutorial)"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "CS309"
$wsA.Range("C7").Value = "CS303"
$wsA.Range("D7").Value = "Free"
$wsA.Range("E7").Value = "CS304"
$wsA.Range("F7").Value = "Free"

$wsA.Range("B8").Value = "Free"
$wsA.Range("C8").Value = "Free"
$wsA.Range("D8").Value = "Free"
$wsA.Range("E8").Value = "CS304 (Tutorial)"
$wsA.Range("F8").Value = "CS303 (Tutorial)"

# --- Section_B (sheet2) ---------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "Free"
$wsB.Range("C2").Value = "CS304"
$wsB.Range("D2").Value = "CS304"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "Free"

$wsB.Range("B3").Value = "CS304"
$wsB.Range("C3").Value = "CS309"
$wsB.Range("D3").Value = "CS303"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "CS309"

$wsB.Range("B5").Value = "CS309"
$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "CS303"
$wsB.Range("F5").Value = "Free"

$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "Free"
$wsB.Range("F6").Value = "CS309 (Tutorial)"

$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "Free"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "CS303"

$wsB.Range("B8").Value = "CS304 (Tutorial)"
$wsB.Range("C8").Value = "Free"
$wsB.Range("D8").Value = "Free"
$wsB.Range("E8").Value = "Free"
$wsB.Range("F8").Value = "CS303 (Tutorial)"
